# Updated cryptos list on Mon Jan 29 06:12:03 UTC 2024 with GitHub Actions
#
# Refreshes the Price / Volume(1h) columns for the existing coin rows, and
# for the rows whose underlying coin changed (46, 47, 49, 50, 51) also
# refreshes the Coin name and Link columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (matching the source data, which is
# always an inline/shared string even when it looks numeric, e.g. "306.67")
# without leaving the cell tagged with a non-default style/number format.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.179.14"
Set-TextValue "E2" "  -0.73%  "
Set-TextValue "D3" "2.266.39"
Set-TextValue "E3" "  -0.99%  "
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "306.67"
Set-TextValue "E5" "  -0.12%  "
Set-TextValue "D6" "96.65"
Set-TextValue "E6" "  +0.74%  "
Set-TextValue "D7" "0.528"
Set-TextValue "E7" "  -0.88%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "E9" "  -1.25%  "
Set-TextValue "D10" "35.05"
Set-TextValue "E10" "  -2.46%  "
Set-TextValue "D11" "0.0791"
Set-TextValue "E11" "  -1.84%  "
Set-TextValue "E12" "  +0.14%  "
Set-TextValue "D13" "6.88"
Set-TextValue "E13" "  +2.13%  "
Set-TextValue "D14" "2.616.74"
Set-TextValue "E14" "  -1.04%  "
Set-TextValue "D15" "14.67"
Set-TextValue "E15" "  +1.10%  "
Set-TextValue "D16" "2.255.09"
Set-TextValue "E16" "  -1.94%  "
Set-TextValue "D17" "0.792"
Set-TextValue "E17" "  -1.44%  "
Set-TextValue "D18" "42.026.84"
Set-TextValue "D19" "12.29"
Set-TextValue "E19" "  -3.43%  "
Set-TextValue "D20" "0.0₃0905"
Set-TextValue "E20" "  -1.74%  "
Set-TextValue "D21" "6.01"
Set-TextValue "E21" "  -0.02%  "
Set-TextValue "D22" "67.80"
Set-TextValue "E22" "  -0.60%  "
Set-TextValue "D23" "237.23"
Set-TextValue "E23" "  -2.43%  "
Set-TextValue "D24" "2.58"
Set-TextValue "E24" "  -1.12%  "
Set-TextValue "E25" "  +0.34%  "
Set-TextValue "E26" "  +0.01%  "
Set-TextValue "D27" "23.51"
Set-TextValue "E27" "  -2.74%  "
Set-TextValue "D28" "37.61"
Set-TextValue "E28" "  +3.95%  "
Set-TextValue "D29" "9.55"
Set-TextValue "E29" "  -0.92%  "
Set-TextValue "E30" "  +0.94%  "
Set-TextValue "D31" "162.44"
Set-TextValue "E31" "  +0.59%  "
Set-TextValue "D32" "5.24"
Set-TextValue "E32" "  -2.08%  "
Set-TextValue "E33" "  +0.01%  "
Set-TextValue "D34" "3.16"
Set-TextValue "E34" "  +2.71%  "
Set-TextValue "D35" "17.69"
Set-TextValue "E35" "  +2.98%  "
Set-TextValue "D36" "0.0737"
Set-TextValue "E36" "  -2.32%  "
Set-TextValue "E38" "  -4.32%  "
Set-TextValue "E39" "  -0.99%  "
Set-TextValue "E40" "  -1.38%  "
Set-TextValue "D41" "4.08"
Set-TextValue "E41" "  -2.95%  "
Set-TextValue "E42" "  +2.93%  "
Set-TextValue "D43" "1.950.73"
Set-TextValue "E43" "  -3.30%  "
Set-TextValue "D44" "18.87"
Set-TextValue "E44" "  -2.91%  "
Set-TextValue "E45" "  -1.09%  "
Set-TextValue "B46" "FraxShare"
Set-TextValue "C46" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D46" "9.92"
Set-TextValue "E46" "  -2.75%  "
Set-TextValue "B47" "NEARProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D47" "2.92"
Set-TextValue "E47" "  -2.73%  "
Set-TextValue "D48" "54.03"
Set-TextValue "E48" "  +1.00%  "
Set-TextValue "B49" "RocketPoolETH"
Set-TextValue "C49" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D49" "2.488.77"
Set-TextValue "E49" "  -0.61%  "
Set-TextValue "B50" "Aave"
Set-TextValue "C50" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "92.07"
Set-TextValue "E50" "  -0.80%  "
Set-TextValue "B51" "BitcoinSV"
Set-TextValue "C51" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D51" "71.58"
Set-TextValue "E51" "  -2.19%  "
